$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear contents of D2:D3 (removes the date value in D2 and the text in D3)
$ws.Range("D2:D3").ClearContents()

# D3 takes on the same number format style as D2 (date format)
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection to D2
$ws.Range("D2").Select()
